# Update "want-to-go" counts (column F) and a few ticket-price cells
# (column G) on the "展览" and "全部类型" sheets, matching the refreshed
# scrape output committed at 456a3b4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 1361
$ws1.Range("F3").Value  = 1292
$ws1.Range("F4").Value  = 905
$ws1.Range("F8").Value  = 116
$ws1.Range("F11").Value = 2452
$ws1.Range("F12").Value = 1608
$ws1.Range("F13").Value = 1494
$ws1.Range("F14").Value = 313
$ws1.Range("G14").Value = "已售罄"
$ws1.Range("F16").Value = 613
$ws1.Range("F17").Value = 786
$ws1.Range("F18").Value = 78
$ws1.Range("F19").Value = 311
$ws1.Range("F23").Value = 529
$ws1.Range("F24").Value = 5022
$ws1.Range("F25").Value = 217
$ws1.Range("F26").Value = 519
$ws1.Range("F27").Value = 78
$ws1.Range("F31").Value = 218
$ws1.Range("F32").Value = 31
$ws1.Range("F34").Value = 737
$ws1.Range("F38").Value = 392
$ws1.Range("F39").Value = 1063
$ws1.Range("F40").Value = 132
$ws1.Range("F42").Value = 173
$ws1.Range("F43").Value = 130
$ws1.Range("F44").Value = 36

# ---------------------------------------------------------------------
# Sheet "全部类型"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value  = 1361
$ws4.Range("F5").Value  = 1292
$ws4.Range("F6").Value  = 905
$ws4.Range("F12").Value = 116
$ws4.Range("F17").Value = 2452
$ws4.Range("F18").Value = 1608
$ws4.Range("F19").Value = 1494
$ws4.Range("F20").Value = 313
$ws4.Range("G20").Value = "已售罄"
$ws4.Range("F22").Value = 613
$ws4.Range("F24").Value = 786
$ws4.Range("F25").Value = 78
$ws4.Range("F26").Value = 311
$ws4.Range("F29").Value = 5022
$ws4.Range("F30").Value = 217
$ws4.Range("F31").Value = 519
$ws4.Range("F32").Value = 78
$ws4.Range("F36").Value = 219
$ws4.Range("F37").Value = 31
$ws4.Range("F39").Value = 737
$ws4.Range("F41").Value = 392
$ws4.Range("F42").Value = 1063
$ws4.Range("F43").Value = 132
$ws4.Range("F44").Value = 173
$ws4.Range("F45").Value = 130
$ws4.Range("F46").Value = 36
